# Applies the 16.2.1 workbook update:
#  - relabels several "by ..." header cells to capitalized "By ..." form
#  - relabels age-related header cells
#  - fills previously blank footnote-reference cells (D28/D29/D31/D32/D33) with "-"
#  - wraps text + increases row height for the two "functional difficulties" header rows (27 and 30)
#  - leaves the final selection on cell B30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------
# (the order below matches the order new shared strings were first introduced
#  in the authored workbook, so the resulting shared-string table lines up)

# Previously-empty footnote cells now show a dash placeholder
$ws.Range("D28").Value = "-"
$ws.Range("D29").Value = "-"
$ws.Range("D31").Value = "-"
$ws.Range("D32").Value = "-"
$ws.Range("D33").Value = "-"

$ws.Range("C6").Value  = "By sex"
$ws.Range("C12").Value = "By territory"
$ws.Range("C22").Value = "By age (in years)"
$ws.Range("C34").Value = "Wealth quintile"
$ws.Range("B22").Value = "По возрасту (в годах)"
$ws.Range("A22").Value = "Жаш курагы боюнча (жылдарда)"

# --- Formatting updates for rows 27 and 30 --------------------------------

$r27 = $ws.Range("A27:B27")
$r27.WrapText = $true
$r27.HorizontalAlignment = -4131
$r27.VerticalAlignment = -4108
$ws.Rows.Item(27).RowHeight = 24

$r30 = $ws.Range("A30:B30")
$r30.WrapText = $true
$r30.HorizontalAlignment = -4131
$r30.VerticalAlignment = -4108
$ws.Rows.Item(30).RowHeight = 36

# --- Final selection --------------------------------------------------------

$ws.Range("B30").Select()
